$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.223.76"
$ws.Range("E2").Value = "  +1.46%  "
$ws.Range("D3").Value = "1.642.71"
$ws.Range("E3").Value = "  +0.32%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("E5").Value = "  +0.16%  "
$ws.Range("E6").Value = "  +1.39%  "
$ws.Range("E8").Value = "  -0.19%  "
$ws.Range("E9").Value = "  +0.82%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.98"
$ws.Range("E10").Value = "  +0.48%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0848"
$ws.Range("E11").Value = "  +0.18%  "
$ws.Range("E12").Value = "  +0.37%  "
$ws.Range("D13").Value = "1.647.02"
$ws.Range("E13").Value = "  +0.49%  "
$ws.Range("E14").Value = "  +0.73%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.543"
$ws.Range("E15").Value = "  +2.66%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "67.10"
$ws.Range("E16").Value = "  +0.47%  "
$ws.Range("D17").Value = "27.221.89"
$ws.Range("E17").Value = "  +1.45%  "
$ws.Range("D18").Value = "0.0₃0741"
$ws.Range("E18").Value = "  +1.93%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "219.15"
$ws.Range("E19").Value = "  +0.02%  "
$ws.Range("E20").Value = "  +0.12%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.98"
$ws.Range("E21").Value = "  +3.84%  "
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.42"
$ws.Range("E22").Value = "  +0.57%  "
$ws.Range("B23").Value = "Toncoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.52"
$ws.Range("E23").Value = "  +3.64%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.13"
$ws.Range("E24").Value = "  -0.28%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "148.05"
$ws.Range("E25").Value = "  +0.69%  "
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.46"
$ws.Range("E27").Value = "  +1.54%  "
$ws.Range("E28").Value = "  -0.34%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.73"
$ws.Range("E29").Value = "  -0.34%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0510"
$ws.Range("E30").Value = "  +1.57%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.18"
$ws.Range("E31").Value = "  +0.67%  "
$ws.Range("E32").Value = "  +1.65%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.01"
$ws.Range("E33").Value = "  +0.36%  "
$ws.Range("D34").Value = "1.308.49"
$ws.Range("E34").Value = "  +3.65%  "
$ws.Range("E35").Value = "  +1.26%  "
$ws.Range("E36").Value = "  +1.41%  "
$ws.Range("E37").Value = "  -0.77%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.552"
$ws.Range("E38").Value = "  +3.67%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.858"
$ws.Range("E39").Value = "  +3.22%  "
$ws.Range("E41").Value = "  +0.71%  "
$ws.Range("E42").Value = "  +6.02%  "
$ws.Range("E43").Value = "  -1.85%  "
$ws.Range("D44").Value = "1.783.99"
$ws.Range("E44").Value = "  +0.12%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "61.88"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "91.87"
$ws.Range("E46").Value = "  -0.17%  "
$ws.Range("E47").Value = "  +1.62%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "0.0₆0107"
$ws.Range("E48").Value = "  +1.42%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0513"
$ws.Range("E49").Value = "  +0.05%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.66"
$ws.Range("E50").Value = "  +0.84%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0966"
$ws.Range("E51").Value = "  +0.60%  "
